# Updated cryptos list - apply Price (D) and Volume(1h) (E) changes per row
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.613.34"
$ws.Range("E2").Value = "  -2.30%  "

$ws.Range("D3").Value = "2.902.85"
$ws.Range("E3").Value = "  -3.28%  "

$ws.Range("E4").Value = "  +0.01%  "

$ws.Range("D5").Value = "'527.86"
$ws.Range("E5").Value = "  -4.22%  "

$ws.Range("D6").Value = "'142.89"
$ws.Range("E6").Value = "  -7.53%  "

$ws.Range("D7").Value = "'1.00"
$ws.Range("E7").Value = "  +0.11%  "

$ws.Range("D8").Value = "'0.555"
$ws.Range("E8").Value = "  -2.05%  "

$ws.Range("D9").Value = "2.908.90"
$ws.Range("E9").Value = "  -3.30%  "

$ws.Range("E10").Value = "  -3.86%  "

$ws.Range("E11").Value = "  -5.33%  "

$ws.Range("E12").Value = "  -2.02%  "

$ws.Range("D13").Value = "3.408.44"
$ws.Range("E13").Value = "  -3.26%  "

$ws.Range("E14").Value = "  +1.64%  "

$ws.Range("D15").Value = "60.582.97"
$ws.Range("E15").Value = "  -2.45%  "

$ws.Range("E16").Value = "  -4.88%  "

$ws.Range("D17").Value = "2.905.59"
$ws.Range("E17").Value = "  -3.17%  "

$ws.Range("E18").Value = "  -4.74%  "

$ws.Range("D19").Value = "'5.02"
$ws.Range("E19").Value = "  -1.96%  "

$ws.Range("E20").Value = "  -2.99%  "

$ws.Range("D21").Value = "'363.10"
$ws.Range("E21").Value = "  -7.74%  "

$ws.Range("E22").Value = "  -0.52%  "

$ws.Range("E23").Value = "  +0.00%  "

$ws.Range("D24").Value = "'63.53"
$ws.Range("E24").Value = "  -2.47%  "

$ws.Range("D25").Value = "3.011.74"
$ws.Range("E25").Value = "  -3.45%  "

$ws.Range("E26").Value = "  -3.46%  "

$ws.Range("E27").Value = "  -1.93%  "

$ws.Range("E28").Value = "  +0.05%  "

$ws.Range("D29").Value = "'7.85"
$ws.Range("E29").Value = "  -7.58%  "

$ws.Range("E30").Value = "  -9.48%  "

$ws.Range("E31").Value = "  -0.02%  "

$ws.Range("D32").Value = "'1.68"
$ws.Range("E32").Value = "  -2.91%  "

$ws.Range("D33").Value = "'19.54"
$ws.Range("E33").Value = "  -4.52%  "

$ws.Range("D34").Value = "'148.91"
$ws.Range("E34").Value = "  -6.75%  "

$ws.Range("E35").Value = "  -6.57%  "

$ws.Range("E36").Value = "  -7.22%  "

$ws.Range("E37").Value = "  -7.49%  "

$ws.Range("E38").Value = "  -6.38%  "

$ws.Range("D39").Value = "'37.91"
$ws.Range("E39").Value = "  +1.90%  "

$ws.Range("E40").Value = "  -4.45%  "

$ws.Range("D41").Value = "2.335.37"
$ws.Range("E41").Value = "  -4.80%  "

$ws.Range("E42").Value = "  -6.32%  "

$ws.Range("D43").Value = "'0.643"
$ws.Range("E43").Value = "  -3.05%  "

$ws.Range("D44").Value = "'20.78"
$ws.Range("E44").Value = "  -7.47%  "

$ws.Range("E45").Value = "  -3.67%  "

$ws.Range("D46").Value = "'0.997"
$ws.Range("E46").Value = "  -0.04%  "

$ws.Range("D47").Value = "'5.04"
$ws.Range("E47").Value = "  +1.65%  "

$ws.Range("E48").Value = "  -4.88%  "

$ws.Range("D49").Value = "'0.0936"
$ws.Range("E49").Value = "  -1.34%  "

$ws.Range("D50").Value = "'10.35"
$ws.Range("E50").Value = "  -1.31%  "

$ws.Range("D51").Value = "'250.71"
$ws.Range("E51").Value = "  -5.97%  "
